# Update the dSF column (F) values on Sheet1 to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 2
    6  = 1
    7  = 7
    9  = -4
    10 = -1
    11 = 6
    12 = -3
    13 = 5
    14 = 2
    16 = -2
    17 = -3
    18 = 4
    19 = 2
    21 = 6
    22 = -6
    24 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
